$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Hyperlinks.Add($ws.Range("C33"), "https://codeforces.com/contest/698/problem/A") | Out-Null
"done"
